$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BPRproductdetails")
$c = $ws.Range("O8")
$f = $c.Font
$f.Name = "Source Sans Pro"
$f.Size = 8
$f.Color = 2236962
Write-Host "done"
